$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Merge "Address Line 1" / "Address Line 2" into a single "Street Address" column.
# Remove the old "Address Line 2" column (D) which shifts everything after it
# one column to the left, preserving widths/formatting of the surviving columns.
$ws.Columns("D").Delete() | Out-Null

# Rename the remaining address column header.
$ws.Range("C1").Value = "Street Address"

# Restore the selection left behind by the edit.
$ws.Range("F5").Select() | Out-Null
